$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 19250
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""

$ws.Range("H23").Value = 19250
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""

$ws.Range("H70").Value = 3624.25
$ws.Range("J70").Value = 3831.3333
$ws.Range("L70").Value = 11493.9999
$ws.Range("N70").Value = -12033.9999

$ws.Range("H73").Value = 3624.25
$ws.Range("J73").Value = 3831.3333
$ws.Range("L73").Value = 11493.9999
$ws.Range("N73").Value = -13365.9999

$ws.Range("H100").Value = 2248.75
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2248.75
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = ""
$ws.Range("M100").Value = 2248.75
$ws.Range("N100").Value = -3330.75

$ws.Range("H113").Value = 12559.571
$ws.Range("I113").Value = 10724.75
$ws.Range("K113").Value = 10724.75
$ws.Range("M113").Value = -7470.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3962.125
$ws.Range("I45").Value = 1961.75
$ws.Range("J45").Value = 5962.5
$ws.Range("K45").Value = 1961.75
$ws.Range("L45").Value = 5962.5
$ws.Range("M45").Value = -1584.75
$ws.Range("N45").Value = -6716.5

$ws.Range("H61").Value = 3449.75
$ws.Range("I61").Value = 3266.3333
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3266.3333
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3054.3333
$ws.Range("N61").Value = -4424

$ws.Range("H136").Value = 3449.75
$ws.Range("I136").Value = 3266.3333
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 9798.999899999999
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -7248.999899999999
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 3182.3333
$ws.Range("I8").Value = 300
$ws.Range("K8").Value = 300
$ws.Range("M8").Value = -160

$ws.Range("H16").Value = 1249.25
$ws.Range("I16").Value = 999.3333
$ws.Range("K16").Value = 999.3333
$ws.Range("M16").Value = -829.3333

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = ""
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = 0

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = ""
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = 0

$ws.Range("H94").Value = 3040.6
$ws.Range("I94").Value = 3012.111
$ws.Range("K94").Value = 3012.111
$ws.Range("M94").Value = -2561.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = ""
$ws.Range("N4").Value = 0

$ws.Range("H9").Value = 229997.5
$ws.Range("J9").Value = 229997.5
$ws.Range("L9").Value = 229997.5
$ws.Range("N9").Value = -230333.5

$ws.Range("H35").Value = 4510.4
$ws.Range("I35").Value = 4381
$ws.Range("K35").Value = 4381
$ws.Range("M35").Value = -4087

$ws.Range("H107").Value = 953
$ws.Range("I107").Value = 793.25
$ws.Range("J107").Value = 1166
$ws.Range("K107").Value = 793.25
$ws.Range("L107").Value = 1166
$ws.Range("M107").Value = 1126.75
$ws.Range("N107").Value = -5006

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 157.16667
$ws.Range("J6").Value = 267
$ws.Range("L6").Value = 801
$ws.Range("N6").Value = -1027

$ws.Range("H26").Value = 360.85715
$ws.Range("J26").Value = 1167.5
$ws.Range("L26").Value = 3502.5
$ws.Range("N26").Value = -4078.5

$ws.Range("H29").Value = 150.5
$ws.Range("I29").Value = 101
$ws.Range("K29").Value = 303
$ws.Range("M29").Value = -26

$ws.Range("H34").Value = 13384.77
$ws.Range("I34").Value = 2
$ws.Range("J34").Value = 14500
$ws.Range("K34").Value = 6
$ws.Range("L34").Value = 43500
$ws.Range("M34").Value = 78
$ws.Range("N34").Value = -43668

$ws.Range("H137").Value = 2049.75
$ws.Range("J137").Value = 2399.6667
$ws.Range("L137").Value = 7199.000100000001
$ws.Range("N137").Value = -17399.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 19997
$ws.Range("J24").Value = 19997
$ws.Range("L24").Value = 19997
$ws.Range("N24").Value = -20343

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = ""
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = 0

$ws.Range("H68").Value = 3599.8333
$ws.Range("J68").Value = 4049.75
$ws.Range("L68").Value = 4049.75
$ws.Range("N68").Value = -5547.75

$ws.Range("H71").Value = 3599.8333
$ws.Range("J71").Value = 4049.75
$ws.Range("L71").Value = 20248.75
$ws.Range("N71").Value = -27736.75

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = ""
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = 0

$ws.Range("H136").Value = 1699999.6
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 5631.8184
$ws.Range("I23").Value = 634.2857
$ws.Range("K23").Value = 634.2857
$ws.Range("M23").Value = -405.2857

$ws.Range("H132").Value = 2416.6667
$ws.Range("I132").Value = 2500
$ws.Range("K132").Value = 7500
$ws.Range("M132").Value = -4970

$ws.Range("H136").Value = 1163.96
$ws.Range("I136").Value = 1083.9048
$ws.Range("J136").Value = 1584.25
$ws.Range("K136").Value = 3251.7144
$ws.Range("L136").Value = 4752.75
$ws.Range("M136").Value = -701.7143999999998
$ws.Range("N136").Value = -9852.75
